# Mixed membership-function test data for FIS rel: swap the Accuracy/Objectivity
# columns (A<->B) for every data row, swap the two header labels to match,
# clear out the now-unused helper column F, and update the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Swap header labels in row 1 (A1 <-> B1) -------------------------------
$h1 = $ws.Cells.Item(1, 1).Value2
$h2 = $ws.Cells.Item(1, 2).Value2
$ws.Cells.Item(1, 1).Value = $h2
$ws.Cells.Item(1, 2).Value = $h1

# --- Swap Accuracy/Objectivity values for every data row (2-27) -----------
for ($r = 2; $r -le 27; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}

# --- Clear the stray/unused column F cells in rows 1-10 --------------------
$ws.Range("F1:F10").Clear()

# --- Update the saved selection to column A ---------------------------------
$ws.Range("A1:A27").Select()

# --- Best-effort: restore the window layout recorded with the edit ---------
$win = $wb.Windows.Item(1)
$win.Left = 945
$win.Top = 11760
$win.Width = 12150
$win.Height = 13410
